$wb = $excel.ActiveWorkbook

# The long "version mismatch" error detail text shared by both locale sheets' Error Detail column.
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c4176dc79c52d29abbc5486d4515699998e92d34/e2e/e86f224a-1064-418f-b39f-10abf8f0b60d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8b6eff97d9f824a1ce146dc0aa27ab5710faee61/e2e/e86f224a-1064-418f-b39f-10abf8f0b60d.md."

$targetFileDisplay = "e86f224a-1064-418f-b39f-10abf8f0b60d.md"
$targetFileUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8b6eff97d9f824a1ce146dc0aa27ab5710faee61/e2e/e86f224a-1064-418f-b39f-10abf8f0b60d.md"

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Widen the "Error Detail" column (P, the 16th column) to fit the new long message.
    $ws.Columns.Item(16).ColumnWidth = 39.17

    # Row 8 (e86f224a-...) just received its handback report:
    #   I8 - Latest Target File: becomes a hyperlink to the handed-back md file.
    $ws.Hyperlinks.Add($ws.Range("I8"), $targetFileUrl, $null, $null, $targetFileDisplay) | Out-Null
    $ws.Range("I8").Font.Underline = $true
    $ws.Range("I8").Font.Color = 15570276

    #   J8 - Latest Handback File: same xlf file name as the latest handoff file (G8).
    $ws.Range("J8").Value = $ws.Range("G8").Value

    #   K8 - Latest Handback DateTime.
    if ($sheetName -eq "zh-cn") {
        $ws.Range("K8").Value = "2016-08-18 06:41:43"
    } else {
        $ws.Range("K8").Value = "2016-08-18 06:41:50"
    }

    #   P8 - Error Detail: handback version is stale compared to latest.
    $ws.Range("P8").Value = $errorDetail
}
